$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated transition-matrix values from games pulled March 7
$updates = @(
    @{ Row = 2; Col = 2; Value = 0.1740890688259109 },
    @{ Row = 2; Col = 3; Value = 0.582995951417004 },
    @{ Row = 2; Col = 10; Value = 0.008097165991902834 },
    @{ Row = 2; Col = 16; Value = 0.1497975708502024 },
    @{ Row = 2; Col = 19; Value = 0.08502024291497975 },
    @{ Row = 3; Col = 2; Value = 0.006535947712418301 },
    @{ Row = 3; Col = 3; Value = 0.0261437908496732 },
    @{ Row = 3; Col = 10; Value = 0.0196078431372549 },
    @{ Row = 3; Col = 16; Value = 0.738562091503268 },
    @{ Row = 3; Col = 19; Value = 0.2091503267973856 },
    @{ Row = 4; Col = 10; Value = 0.02564102564102564 },
    @{ Row = 4; Col = 16; Value = 0.7948717948717948 },
    @{ Row = 4; Col = 19; Value = 0.1794871794871795 },
    @{ Row = 5; Col = 10; Value = 0.3333333333333333 },
    @{ Row = 5; Col = 16; Value = 0.3333333333333333 },
    @{ Row = 5; Col = 19; Value = 0.3333333333333333 },
    @{ Row = 6; Col = 2; Value = 0.0653061224489796 },
    @{ Row = 6; Col = 4; Value = 0.00816326530612245 },
    @{ Row = 6; Col = 6; Value = 0.08571428571428572 },
    @{ Row = 6; Col = 10; Value = 0.2081632653061224 },
    @{ Row = 6; Col = 15; Value = 0.02448979591836735 },
    @{ Row = 6; Col = 17; Value = 0.1591836734693877 },
    @{ Row = 6; Col = 18; Value = 0.08571428571428572 },
    @{ Row = 6; Col = 19; Value = 0.363265306122449 },
    @{ Row = 7; Col = 2; Value = 0.1258278145695364 },
    @{ Row = 7; Col = 4; Value = 0.02649006622516556 },
    @{ Row = 7; Col = 6; Value = 0.06622516556291391 },
    @{ Row = 7; Col = 10; Value = 0.1059602649006623 },
    @{ Row = 7; Col = 15; Value = 0.02649006622516556 },
    @{ Row = 7; Col = 17; Value = 0.1986754966887417 },
    @{ Row = 7; Col = 18; Value = 0.0728476821192053 },
    @{ Row = 7; Col = 19; Value = 0.3774834437086093 },
    @{ Row = 8; Col = 2; Value = 0.08057851239669421 },
    @{ Row = 8; Col = 4; Value = 0.02272727272727273 },
    @{ Row = 8; Col = 6; Value = 0.06611570247933884 },
    @{ Row = 8; Col = 10; Value = 0.08677685950413223 },
    @{ Row = 8; Col = 15; Value = 0.01239669421487603 },
    @{ Row = 8; Col = 17; Value = 0.1921487603305785 },
    @{ Row = 8; Col = 18; Value = 0.1053719008264463 },
    @{ Row = 8; Col = 19; Value = 0.4338842975206612 },
    @{ Row = 9; Col = 2; Value = 0.08482142857142858 },
    @{ Row = 9; Col = 4; Value = 0.01785714285714286 },
    @{ Row = 9; Col = 6; Value = 0.07142857142857142 },
    @{ Row = 9; Col = 10; Value = 0.06696428571428571 },
    @{ Row = 9; Col = 15; Value = 0.01785714285714286 },
    @{ Row = 9; Col = 17; Value = 0.1964285714285714 },
    @{ Row = 9; Col = 18; Value = 0.1071428571428571 },
    @{ Row = 9; Col = 19; Value = 0.4375 },
    @{ Row = 10; Col = 2; Value = 0.08825978351373855 },
    @{ Row = 10; Col = 4; Value = 0.01582014987510408 },
    @{ Row = 10; Col = 5; Value = 0.003330557868442964 },
    @{ Row = 10; Col = 6; Value = 0.07077435470441298 },
    @{ Row = 10; Col = 10; Value = 0.09492089925062448 },
    @{ Row = 10; Col = 15; Value = 0.01082431307243963 },
    @{ Row = 10; Col = 17; Value = 0.2081598667776852 },
    @{ Row = 10; Col = 18; Value = 0.1074104912572856 },
    @{ Row = 10; Col = 19; Value = 0.4004995836802664 },
    @{ Row = 11; Col = 7; Value = 0.1363636363636364 },
    @{ Row = 11; Col = 10; Value = 0.1212121212121212 },
    @{ Row = 11; Col = 11; Value = 0.2348484848484849 },
    @{ Row = 11; Col = 12; Value = 0.4962121212121212 },
    @{ Row = 11; Col = 19; Value = 0.01136363636363636 },
    @{ Row = 12; Col = 7; Value = 0.7299270072992701 },
    @{ Row = 12; Col = 10; Value = 0.2262773722627737 },
    @{ Row = 12; Col = 12; Value = 0.0291970802919708 },
    @{ Row = 12; Col = 19; Value = 0.0145985401459854 },
    @{ Row = 13; Col = 7; Value = 0.75 },
    @{ Row = 13; Col = 10; Value = 0.21875 },
    @{ Row = 13; Col = 19; Value = 0.03125 },
    @{ Row = 15; Col = 6; Value = 0.01809954751131222 },
    @{ Row = 15; Col = 8; Value = 0.1538461538461539 },
    @{ Row = 15; Col = 9; Value = 0.08144796380090498 },
    @{ Row = 15; Col = 10; Value = 0.4072398190045249 },
    @{ Row = 15; Col = 11; Value = 0.05429864253393665 },
    @{ Row = 15; Col = 13; Value = 0.009049773755656109 },
    @{ Row = 15; Col = 15; Value = 0.09954751131221719 },
    @{ Row = 15; Col = 19; Value = 0.1764705882352941 },
    @{ Row = 16; Col = 6; Value = 0.03932584269662921 },
    @{ Row = 16; Col = 8; Value = 0.2022471910112359 },
    @{ Row = 16; Col = 9; Value = 0.101123595505618 },
    @{ Row = 16; Col = 10; Value = 0.3539325842696629 },
    @{ Row = 16; Col = 11; Value = 0.1348314606741573 },
    @{ Row = 16; Col = 13; Value = 0.03370786516853932 },
    @{ Row = 16; Col = 14; Value = 0.005617977528089887 },
    @{ Row = 16; Col = 15; Value = 0.03370786516853932 },
    @{ Row = 16; Col = 19; Value = 0.09550561797752809 },
    @{ Row = 17; Col = 6; Value = 0.0155902004454343 },
    @{ Row = 17; Col = 8; Value = 0.1937639198218263 },
    @{ Row = 17; Col = 9; Value = 0.09354120267260579 },
    @{ Row = 17; Col = 10; Value = 0.4075723830734966 },
    @{ Row = 17; Col = 11; Value = 0.06904231625835189 },
    @{ Row = 17; Col = 13; Value = 0.0111358574610245 },
    @{ Row = 17; Col = 15; Value = 0.08463251670378619 },
    @{ Row = 17; Col = 19; Value = 0.1247216035634744 },
    @{ Row = 18; Col = 6; Value = 0.02941176470588235 },
    @{ Row = 18; Col = 8; Value = 0.2142857142857143 },
    @{ Row = 18; Col = 9; Value = 0.08403361344537816 },
    @{ Row = 18; Col = 10; Value = 0.4285714285714285 },
    @{ Row = 18; Col = 11; Value = 0.08823529411764706 },
    @{ Row = 18; Col = 13; Value = 0.01260504201680672 },
    @{ Row = 18; Col = 15; Value = 0.05882352941176471 },
    @{ Row = 18; Col = 19; Value = 0.08403361344537816 },
    @{ Row = 19; Col = 6; Value = 0.02365930599369085 },
    @{ Row = 19; Col = 8; Value = 0.2200315457413249 },
    @{ Row = 19; Col = 9; Value = 0.09779179810725552 },
    @{ Row = 19; Col = 10; Value = 0.3682965299684542 },
    @{ Row = 19; Col = 11; Value = 0.08753943217665615 },
    @{ Row = 19; Col = 13; Value = 0.01419558359621451 },
    @{ Row = 19; Col = 15; Value = 0.06861198738170347 },
    @{ Row = 19; Col = 19; Value = 0.1198738170347003 }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, $u.Col).Value = $u.Value
}

Write-Output "Applied $($updates.Count) cell updates"
